$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.039.03'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").Value = '2.311.97'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.69'
$ws.Range("E5").Value = '  +4.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.62'
$ws.Range("E6").Value = '  +3.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.534'
$ws.Range("E8").Value = '  +1.69%  '
$ws.Range("D9").Value = '2.335.37'
$ws.Range("E9").Value = '  +2.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("E10").Value = '  +8.21%  '
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.15'
$ws.Range("E12").Value = '  +7.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.342'
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.99'
$ws.Range("E14").Value = '  +5.03%  '
$ws.Range("D15").Value = '2.725.33'
$ws.Range("E15").Value = '  +2.19%  '
$ws.Range("D16").Value = '56.255.30'
$ws.Range("E16").Value = '  +3.65%  '
$ws.Range("E17").Value = '  +4.39%  '
$ws.Range("D18").Value = '2.344.37'
$ws.Range("E18").Value = '  +3.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.50'
$ws.Range("E19").Value = '  +2.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.28'
$ws.Range("E20").Value = '  +3.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.88'
$ws.Range("E21").Value = '  +6.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.65'
$ws.Range("E22").Value = '  +5.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.01'
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.995'
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  +6.70%  '
$ws.Range("E27").Value = '  +4.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.25'
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.72'
$ws.Range("E29").Value = '  +7.56%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.20'
$ws.Range("E30").Value = '  +10.53%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0727'
$ws.Range("E31").Value = '  +5.28%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.27'
$ws.Range("E32").Value = '  +5.37%  '
$ws.Range("E33").Value = '  +3.10%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.993'
$ws.Range("E35").Value = '  -0.44%  '
$ws.Range("E36").Value = '  +5.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.926'
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("E38").Value = '  +8.03%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.39'
$ws.Range("E39").Value = '  +4.03%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.51'
$ws.Range("E40").Value = '  +8.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.383'
$ws.Range("E41").Value = '  +2.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '139.39'
$ws.Range("E42").Value = '  +11.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.58'
$ws.Range("E43").Value = '  +6.31%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '271.98'
$ws.Range("E44").Value = '  +13.05%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.03'
$ws.Range("E45").Value = '  +4.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0509'
$ws.Range("E46").Value = '  +3.27%  '
$ws.Range("E47").Value = '  +3.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.554'
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("E49").Value = '  +2.40%  '
$ws.Range("E50").Value = '  +5.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.88'
$ws.Range("E51").Value = '  +4.80%  '
